# Generate Report for Handoff
#
# Rows 4-7 in the zh-cn and de-de sheets represent files that were still
# "Ready for handoff" (Priority = low). Running the handoff-report
# generation process for them:
#   - bumps their Priority from "low" to "ht"
#   - stamps a fresh "Latest Handoff Datetime" / "Latest HO Xliff Generate
#     Date" on each of them (the de-de handoff datetime is shared text
#     with the Overview sheet's generate-date column, so that one updates
#     too)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = 4, 5, 6, 7

foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-12 18:36:42"

    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-12 18:36:49"

    $overview.Cells.Item($r, 7).Value = "2016-08-12 18:36:49"
}
